# Apply the updated transition probability matrix values
# (added more games / sped up simulate-game logic / drafted optimization logic)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.55
$ws.Range("P2").Value = 0.15
$ws.Range("S2").Value = 0.1
$ws.Range("J3").Value = 0.1818181818181818
$ws.Range("P3").Value = 0.5454545454545454
$ws.Range("S3").Value = 0.2727272727272727
$ws.Range("F6").Value = 0.1666666666666667
$ws.Range("J6").Value = 0.3333333333333333
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("J7").Value = 0.6666666666666666
$ws.Range("Q7").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.25
$ws.Range("F8").Value = 0.08333333333333333
$ws.Range("J8").Value = 0.1666666666666667
$ws.Range("Q8").Value = 0.3333333333333333
$ws.Range("S8").Value = 0.1666666666666667
$ws.Range("B9").Value = 0.1111111111111111
$ws.Range("F9").Value = 0.1111111111111111
$ws.Range("J9").Value = 0.3333333333333333
$ws.Range("Q9").Value = 0.2222222222222222
$ws.Range("S9").Value = 0.2222222222222222
$ws.Range("B10").Value = 0.1641791044776119
$ws.Range("F10").Value = 0.01492537313432836
$ws.Range("J10").Value = 0.1940298507462687
$ws.Range("O10").Value = 0.01492537313432836
$ws.Range("Q10").Value = 0.2985074626865671
$ws.Range("R10").Value = 0.1044776119402985
$ws.Range("S10").Value = 0.208955223880597
$ws.Range("G11").Value = 0.2222222222222222
$ws.Range("J11").Value = 0.2222222222222222
$ws.Range("K11").Value = 0.4444444444444444
$ws.Range("L11").Value = 0.1111111111111111
$ws.Range("G12").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.1666666666666667
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("O15").Value = 0.3333333333333333
$ws.Range("H16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.6666666666666666
$ws.Range("S16").Value = 0.2222222222222222
$ws.Range("H17").Value = 0.1851851851851852
$ws.Range("I17").Value = 0.1481481481481481
$ws.Range("J17").Value = 0.4444444444444444
$ws.Range("K17").Value = 0.07407407407407407
$ws.Range("S17").Value = 0.1481481481481481
$ws.Range("J18").Value = 0.5714285714285714
$ws.Range("K18").Value = 0.2857142857142857
$ws.Range("O18").Value = 0.1428571428571428
$ws.Range("H19").Value = 0.15625
$ws.Range("I19").Value = 0.125
$ws.Range("J19").Value = 0.53125
$ws.Range("K19").Value = 0.03125
$ws.Range("M19").Value = 0.03125
$ws.Range("O19").Value = 0.0625
$ws.Range("S19").Value = 0.0625
